$d = $word.ActiveDocument

# --- Locate the target paragraph ("Gioco adattato ... Expedition.") robustly ---
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Gioco adattato*") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Target paragraph not found"
}

$para = $d.Paragraphs.Item($targetIndex)
$pRange = $para.Range
$paraStart = $pRange.Start
$paraEnd = $pRange.End

# Inner range excludes the trailing paragraph mark (pilcrow)
$inner = $d.Range($paraStart, $paraEnd - 1)

$newRunsXml = '<w:r><w:rPr><w:rFonts w:ascii="Futura Bk BT" w:hAnsi="Futura Bk BT"/><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Gioco basato su 1572: The Lost </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Futura Bk BT" w:hAnsi="Futura Bk BT"/><w:i/><w:iCs/></w:rPr><w:t>Expedition</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Futura Bk BT" w:hAnsi="Futura Bk BT"/><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">, per gentile concessione di Mike </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Futura Bk BT" w:hAnsi="Futura Bk BT"/><w:i/><w:iCs/></w:rPr><w:t>Heiman</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Futura Bk BT" w:hAnsi="Futura Bk BT"/><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Futura Bk BT" w:hAnsi="Futura Bk BT"/><w:i/><w:iCs/></w:rPr><w:t>kzinti</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Futura Bk BT" w:hAnsi="Futura Bk BT"/><w:i/><w:iCs/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Futura Bk BT" w:hAnsi="Futura Bk BT"/><w:i/><w:iCs/></w:rPr><w:br/><w:t>Adattamento di Andrea Gottardi (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Futura Bk BT" w:hAnsi="Futura Bk BT"/><w:i/><w:iCs/></w:rPr><w:t>AndreaGot</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Futura Bk BT" w:hAnsi="Futura Bk BT"/><w:i/><w:iCs/></w:rPr><w:t>) – Mattia Ferrari (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Futura Bk BT" w:hAnsi="Futura Bk BT"/><w:i/><w:iCs/></w:rPr><w:t>mattocrazy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Futura Bk BT" w:hAnsi="Futura Bk BT"/><w:i/><w:iCs/></w:rPr><w:t>)- Giulia Mezzasalma (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Futura Bk BT" w:hAnsi="Futura Bk BT"/><w:i/><w:iCs/></w:rPr><w:t>GiuliaJuliette</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Futura Bk BT" w:hAnsi="Futura Bk BT"/><w:i/><w:iCs/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Futura Bk BT" w:hAnsi="Futura Bk BT"/><w:i/><w:iCs/></w:rPr><w:br/><w:t>In collaborazione con MUSE – Museo delle Scienze di Trento</w:t></w:r>'

$pkgXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $newRunsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$inner.InsertXML($pkgXml)

# --- Re-locate the (now-updated) paragraph and remove the following empty paragraph ---
$para2 = $d.Paragraphs.Item($targetIndex)
$nextPara = $para2.Next()
$nextRange = $nextPara.Range
$nextTextTrimmed = $nextRange.Text.Trim("`r", "`a")
if ($nextTextTrimmed -eq "") {
    $nextRange.Delete()
}

# --- Apply sz=20 / szCs=20 across the whole paragraph (pPr + all runs) ---
$para3 = $d.Paragraphs.Item($targetIndex)
$finalRange = $para3.Range
$finalRange.Font.Size = 10
$finalRange.Font.SizeBi = 10

Write-Output "Done. Final paragraph text:"
Write-Output $finalRange.Text
